$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.779999999999999
$ws.Range("C7").Value = -13.57
$ws.Range("B10").Value = 5.599
$ws.Range("B12").Value = 5.86
$ws.Range("C15").Value = -13.741
$ws.Range("B18").Value = 5.249
$ws.Range("D18").Value = -8.784000000000001
$ws.Range("D19").Value = -8.113000000000001
$ws.Range("C20").Value = -12.673
$ws.Range("D27").Value = -8.745000000000001
$ws.Range("C29").Value = -12.107
$ws.Range("C30").Value = -12.981
$ws.Range("C31").Value = -13.024
$ws.Range("B37").Value = 8.494
$ws.Range("C40").Value = -12.782
$ws.Range("D42").Value = -8.468
$ws.Range("D44").Value = -7.911
$ws.Range("D47").Value = -7.408999999999999
$ws.Range("B55").Value = 5.194
$ws.Range("D58").Value = -8.403
$ws.Range("B68").Value = 5.403
$ws.Range("C68").Value = -11.083
$ws.Range("D73").Value = -7.841999999999999
$ws.Range("C76").Value = -12.977
$ws.Range("B77").Value = 6.235000000000001
$ws.Range("B78").Value = 7.631
$ws.Range("C87").Value = -13.393
$ws.Range("C88").Value = -13.411
$ws.Range("D95").Value = -7.719999999999999
$ws.Range("C96").Value = -12.628
$ws.Range("C98").Value = -13.23
$ws.Range("C101").Value = -12.612
$ws.Range("D101").Value = -7.768000000000001
$ws.Range("C102").Value = -13.086
